# Regenerate the "K" column (column G) values in the save_data sheet.
# The K column used to hold a different metric ("Strike#"); it has been
# recalculated/re-pulled and written back with new values (rows 25 & 26
# already held the correct values and are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 3
    23 = 2
    24 = 2
    27 = 1
    28 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
